# Auto-generated edit script: updates cryptos price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "27.521.48"
Set-TextValue $ws.Range("E2") "  -0.58%  "
Set-TextValue $ws.Range("D3") "1.837.01"
Set-TextValue $ws.Range("E3") "  -0.58%  "
Set-TextValue $ws.Range("D4") "1.008"
Set-TextValue $ws.Range("E4") "  -2.54%  "
Set-TextValue $ws.Range("D5") "316.20"
Set-TextValue $ws.Range("E5") "  -1.57%  "
Set-TextValue $ws.Range("D6") "1.006"
Set-TextValue $ws.Range("E6") "  -2.33%  "
Set-TextValue $ws.Range("D7") "0.4294"
Set-TextValue $ws.Range("D8") "0.3716"
Set-TextValue $ws.Range("D9") "0.07278"
Set-TextValue $ws.Range("E9") "  -1.42%  "
Set-TextValue $ws.Range("D10") "0.8683"
Set-TextValue $ws.Range("E10") "  -1.46%  "
Set-TextValue $ws.Range("D11") "21.23"
Set-TextValue $ws.Range("E11") "  -1.25%  "
Set-TextValue $ws.Range("D12") "1.863.78"
Set-TextValue $ws.Range("E12") "  -0.63%  "
Set-TextValue $ws.Range("D13") "6.704"
Set-TextValue $ws.Range("E13") "  +0.33%  "
Set-TextValue $ws.Range("D14") "5.367"
Set-TextValue $ws.Range("E14") "  -2.22%  "
Set-TextValue $ws.Range("D15") "0.07108"
Set-TextValue $ws.Range("E15") "  -0.22%  "
Set-TextValue $ws.Range("D16") "88.63"
Set-TextValue $ws.Range("E16") "  +4.56%  "
Set-TextValue $ws.Range("E17") "  -2.68%  "
Set-TextValue $ws.Range("D18") "0.000008950"
Set-TextValue $ws.Range("E18") "  -1.26%  "
Set-TextValue $ws.Range("E19") "  -2.47%  "
Set-TextValue $ws.Range("D20") "15.28"
Set-TextValue $ws.Range("E20") "  -1.02%  "
Set-TextValue $ws.Range("D21") "27.530.02"
Set-TextValue $ws.Range("E21") "  -0.63%  "
Set-TextValue $ws.Range("D22") "5.173"
Set-TextValue $ws.Range("E22") "  -1.97%  "
Set-TextValue $ws.Range("D23") "10.96"
Set-TextValue $ws.Range("E23") "  -2.61%  "
Set-TextValue $ws.Range("D24") "2.066.45"
Set-TextValue $ws.Range("E24") "  -1.09%  "
Set-TextValue $ws.Range("D25") "2.008"
Set-TextValue $ws.Range("E25") "  -2.16%  "
Set-TextValue $ws.Range("D26") "153.95"
Set-TextValue $ws.Range("E26") "  -2.89%  "
Set-TextValue $ws.Range("D27") "18.50"
Set-TextValue $ws.Range("E27") "  -0.96%  "
Set-TextValue $ws.Range("D28") "2.167"
Set-TextValue $ws.Range("E28") "  +8.99%  "
Set-TextValue $ws.Range("D29") "5.301"
Set-TextValue $ws.Range("E29") "  -0.41%  "
Set-TextValue $ws.Range("D30") "117.28"
Set-TextValue $ws.Range("E30") "  -0.34%  "
Set-TextValue $ws.Range("D31") "0.08878"
Set-TextValue $ws.Range("E31") "  -2.09%  "
Set-TextValue $ws.Range("D32") "1.211"
Set-TextValue $ws.Range("E32") "  +0.50%  "
Set-TextValue $ws.Range("D33") "0.7702"
Set-TextValue $ws.Range("E33") "  -0.13%  "
Set-TextValue $ws.Range("D34") "4.499"
Set-TextValue $ws.Range("E34") "  -1.06%  "
Set-TextValue $ws.Range("D35") "2.903"
Set-TextValue $ws.Range("E35") "  -3.12%  "
Set-TextValue $ws.Range("E36") "  -2.49%  "
Set-TextValue $ws.Range("E37") "  -2.20%  "
Set-TextValue $ws.Range("D38") "0.05295"
Set-TextValue $ws.Range("E38") "  +0.82%  "
Set-TextValue $ws.Range("D39") "0.01961"
Set-TextValue $ws.Range("E39") "  -0.47%  "
Set-TextValue $ws.Range("D40") "7.158"
Set-TextValue $ws.Range("E40") "  +4.20%  "
Set-TextValue $ws.Range("D41") "2.883"
Set-TextValue $ws.Range("E41") "  +1.15%  "
Set-TextValue $ws.Range("B42") "TheSandbox"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D42") "0.5099"
Set-TextValue $ws.Range("E42") "  -1.30%  "
Set-TextValue $ws.Range("B43") "Algorand"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D43") "0.1679"
Set-TextValue $ws.Range("E43") "  +0.76%  "
Set-TextValue $ws.Range("D44") "8.700"
Set-TextValue $ws.Range("E44") "  +0.09%  "
Set-TextValue $ws.Range("D45") "10.60"
Set-TextValue $ws.Range("E45") "  -1.33%  "
Set-TextValue $ws.Range("D46") "106.69"
Set-TextValue $ws.Range("E46") "  -2.97%  "
Set-TextValue $ws.Range("D47") "0.4727"
Set-TextValue $ws.Range("E47") "  +0.80%  "
Set-TextValue $ws.Range("D49") "1.006"
Set-TextValue $ws.Range("E49") "  -2.46%  "
Set-TextValue $ws.Range("D50") "1.678"
Set-TextValue $ws.Range("E50") "  -1.26%  "
Set-TextValue $ws.Range("D51") "1.835"
Set-TextValue $ws.Range("E51") "  -2.33%  "
